$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 282 (shifts existing rows 282-310 down to 283-311)
$ws.Rows.Item(282).Insert()

# Fill the newly inserted row 282 with the new weekly record
$ws.Cells.Item(282, 1).Value = 11
$ws.Cells.Item(282, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(282, 3).Value = "Bíobío"
$ws.Cells.Item(282, 4).Value = 45212
$ws.Cells.Item(282, 5).Value = 8
$ws.Cells.Item(282, 6).Value = "Fruta"
$ws.Cells.Item(282, 7).Value = 100108
$ws.Cells.Item(282, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(282, 9).Value = 100108005
$ws.Cells.Item(282, 10).Value = "Piña"
$ws.Cells.Item(282, 11).Value = "Caramelo"
$ws.Cells.Item(282, 12).Value = "Segunda"
$ws.Cells.Item(282, 13).Value = 100
$ws.Cells.Item(282, 14).Value = 22000
$ws.Cells.Item(282, 15).Value = 23000
$ws.Cells.Item(282, 16).Value = 22500
$ws.Cells.Item(282, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(282, 18).Value = "Ecuador"
$ws.Cells.Item(282, 19).Value = 1607
$ws.Cells.Item(282, 20).Value = 14
